$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row ("Tanggal" + day-of-month numbers). This shifts every
# subsequent row up by one (so each row now shows what used to be in the row
# below it) and drops the final now-superfluous row from the bottom of the
# used range.
$ws.Rows(1).Delete()

# Update the saved cursor/selection position to match the author's final view.
$null = $ws.Range("E9").Select()
